$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''36.565.79'
$ws.Cells.Item(2, 5).Value = '  -2.20%  '
$ws.Cells.Item(3, 4).Value = '''1.993.63'
$ws.Cells.Item(3, 5).Value = '  -1.15%  '
$ws.Cells.Item(4, 5).Value = '  +0.12%  '
$ws.Cells.Item(5, 4).Value = '''235.79'
$ws.Cells.Item(5, 5).Value = '  -9.46%  '
$ws.Cells.Item(6, 4).Value = '''0.599'
$ws.Cells.Item(6, 5).Value = '  -2.82%  '
$ws.Cells.Item(7, 5).Value = '  +0.09%  '
$ws.Cells.Item(8, 4).Value = '''54.82'
$ws.Cells.Item(8, 5).Value = '  -2.98%  '
$ws.Cells.Item(9, 5).Value = '  -4.23%  '
$ws.Cells.Item(10, 4).Value = '''57.91'
$ws.Cells.Item(10, 5).Value = '  +2.17%  '
$ws.Cells.Item(11, 4).Value = '''0.0747'
$ws.Cells.Item(11, 5).Value = '  -3.54%  '
$ws.Cells.Item(12, 5).Value = '  -3.02%  '
$ws.Cells.Item(13, 5).Value = '  -1.03%  '
$ws.Cells.Item(14, 4).Value = '''2.288.38'
$ws.Cells.Item(14, 5).Value = '  -1.15%  '
$ws.Cells.Item(15, 4).Value = '''20.40'
$ws.Cells.Item(15, 5).Value = '  -3.14%  '
$ws.Cells.Item(16, 5).Value = '  -6.10%  '
$ws.Cells.Item(17, 4).Value = '''5.08'
$ws.Cells.Item(17, 5).Value = '  -3.72%  '
$ws.Cells.Item(18, 4).Value = '''1.995.17'
$ws.Cells.Item(18, 5).Value = '  -1.40%  '
$ws.Cells.Item(19, 4).Value = '''36.512.25'
$ws.Cells.Item(19, 5).Value = '  -2.29%  '
$ws.Cells.Item(20, 4).Value = '''67.84'
$ws.Cells.Item(20, 5).Value = '  -3.06%  '
$ws.Cells.Item(21, 4).Value = '''0.0₃0805'
$ws.Cells.Item(21, 5).Value = '  -4.28%  '
$ws.Cells.Item(22, 4).Value = '''5.27'
$ws.Cells.Item(22, 5).Value = '  +1.29%  '
$ws.Cells.Item(23, 4).Value = '''221.86'
$ws.Cells.Item(23, 5).Value = '  -3.20%  '
$ws.Cells.Item(24, 4).Value = '''1.00'
$ws.Cells.Item(24, 5).Value = '  -0.03%  '
$ws.Cells.Item(25, 4).Value = '''2.37'
$ws.Cells.Item(25, 5).Value = '  +0.36%  '
$ws.Cells.Item(26, 4).Value = '''2.39'
$ws.Cells.Item(26, 5).Value = '  -9.85%  '
$ws.Cells.Item(27, 4).Value = '''162.34'
$ws.Cells.Item(27, 5).Value = '  -1.64%  '
$ws.Cells.Item(28, 4).Value = '''8.66'
$ws.Cells.Item(28, 5).Value = '  -4.03%  '
$ws.Cells.Item(29, 4).Value = '''0.128'
$ws.Cells.Item(29, 5).Value = '  -2.55%  '
$ws.Cells.Item(30, 4).Value = '''18.88'
$ws.Cells.Item(30, 5).Value = '  -5.54%  '
$ws.Cells.Item(31, 5).Value = '  +0.54%  '
$ws.Cells.Item(32, 4).Value = '''0.116'
$ws.Cells.Item(32, 5).Value = '  -3.26%  '
$ws.Cells.Item(33, 4).Value = '''4.37'
$ws.Cells.Item(33, 5).Value = '  -6.62%  '
$ws.Cells.Item(34, 4).Value = '''0.0606'
$ws.Cells.Item(34, 5).Value = '  -6.72%  '
$ws.Cells.Item(35, 4).Value = '''4.25'
$ws.Cells.Item(35, 5).Value = '  -7.35%  '
$ws.Cells.Item(36, 4).Value = '''2.33'
$ws.Cells.Item(36, 5).Value = '  -2.91%  '
$ws.Cells.Item(37, 2).Value = 'BinanceUSD'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(37, 4).Value = '''1.00'
$ws.Cells.Item(37, 5).Value = '  +0.12%  '
$ws.Cells.Item(38, 2).Value = 'RenderToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(38, 4).Value = '''3.37'
$ws.Cells.Item(38, 5).Value = '  +0.18%  '
$ws.Cells.Item(39, 5).Value = '  -2.93%  '
$ws.Cells.Item(40, 4).Value = '''5.61'
$ws.Cells.Item(40, 5).Value = '  +6.52%  '
$ws.Cells.Item(41, 5).Value = '  -1.69%  '
$ws.Cells.Item(42, 5).Value = '  +1.02%  '
$ws.Cells.Item(43, 4).Value = '''1.454.37'
$ws.Cells.Item(43, 5).Value = '  +3.50%  '
$ws.Cells.Item(44, 4).Value = '''0.0202'
$ws.Cells.Item(44, 5).Value = '  -5.26%  '
$ws.Cells.Item(45, 5).Value = '  -8.57%  '
$ws.Cells.Item(46, 4).Value = '''89.10'
$ws.Cells.Item(46, 5).Value = '  -1.39%  '
$ws.Cells.Item(47, 4).Value = '''15.24'
$ws.Cells.Item(47, 5).Value = '  -3.69%  '
$ws.Cells.Item(48, 4).Value = '''0.994'
$ws.Cells.Item(48, 5).Value = '  -3.43%  '
$ws.Cells.Item(49, 5).Value = '  -0.71%  '
$ws.Cells.Item(50, 4).Value = '''6.86'
$ws.Cells.Item(50, 5).Value = '  -3.83%  '
$ws.Cells.Item(51, 4).Value = '''3.72'
$ws.Cells.Item(51, 5).Value = '  +7.82%  '
